$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the selection on the existing "booking_data" sheet: the
#    user clicked cell C9 and then selected the whole sheet (Ctrl+A),
#    so the stored selection becomes the full grid while the scroll/
#    anchor position stays at C9.
# ---------------------------------------------------------------------
$bookingData = $wb.Worksheets.Item("booking_data")
$bookingData.Activate()
$bookingData.Range("C9").Select()
$bookingData.Cells.Select()

# ---------------------------------------------------------------------
# 2) Add the new "Book" worksheet after "booking_data" (becomes the
#    last / active tab).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Book"

# Reuse the exact cell formatting (border/fill/font/number-format)
# already defined in "booking_data" so no new style entries are
# created and the look matches the rest of the workbook.
$bookingData.Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)   # header row style

$bookingData.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)   # plain bordered numeric style

$bookingData.Range("D2").Copy()
$newSheet.Range("D2:D3").PasteSpecial(-4122)   # text ("@") style
$newSheet.Range("B2:C3").PasteSpecial(-4122)   # text ("@") style

$bookingData.Range("E2").Copy()
$newSheet.Range("E2:F3").PasteSpecial(-4122)   # date style

# Column widths (to match the look of the pasted book data).
$newSheet.Columns.Item(2).ColumnWidth = 10.92
$newSheet.Columns.Item(3).ColumnWidth = 16.59
$newSheet.Columns.Item(4).ColumnWidth = 9.25
$newSheet.Columns.Item(5).ColumnWidth = 9.92
$newSheet.Columns.Item(6).ColumnWidth = 22.25

# Header row.
$newSheet.Range("A1").Value = "id"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "description"
$newSheet.Range("D1").Value = "pageCount"
$newSheet.Range("E1").Value = "excerpt"
$newSheet.Range("F1").Value = "publishDate"

# Row 2 - "The Sign of Four".
$newSheet.Range("A2").Value = 983242
$newSheet.Range("F2").Value = "2023-01-02T16:51:12.898Z"
$newSheet.Range("B2").Value = "Sherlok"
$newSheet.Range("C2").Value = "The Sign of Four"
$newSheet.Range("D2").Value = "300"
$newSheet.Range("E2").Value = "Thriller"

# Row 3 - "The Dangerous Path".
$newSheet.Range("A3").Value = 879374
$newSheet.Range("B3").Value = "Romana"
$newSheet.Range("C3").Value = "The Dangerous Path"
$newSheet.Range("D3").Value = "400"
$newSheet.Range("E3").Value = "SiFi"
$newSheet.Range("F3").Value = "2023-03-03T16:51:12.898Z"

# Leave the cursor on E11, as in the authored workbook.
$newSheet.Range("E11").Select()
